$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text (string) representation, matching the source data
$cells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'E12', 'D13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'E22', 'D23', 'B24', 'C24', 'D24', 'E24', 'B25', 'C25', 'D25', 'E25', 'B26', 'C26', 'D26', 'E26', 'B27', 'C27', 'D27', 'E27', 'B28', 'C28', 'D28', 'E28', 'B29', 'C29', 'D29', 'E29', 'B30', 'C30', 'D30', 'E30', 'B31', 'C31', 'D31', 'E31', 'B32', 'C32', 'D32', 'E32', 'B33', 'C33', 'D33', 'E33', 'B34', 'C34', 'D34', 'E34', 'B35', 'C35', 'D35', 'E35', 'B36', 'C36', 'D36', 'E36', 'B37', 'C37', 'D37', 'E37', 'B38', 'C38', 'D38', 'E38', 'B39', 'C39', 'D39', 'E39', 'B40', 'C40', 'D40', 'E40', 'B41', 'C41', 'D41', 'E41', 'B42', 'C42', 'D42', 'E42', 'B43', 'C43', 'D43', 'E43', 'B44', 'C44', 'D44', 'E44', 'B45', 'C45', 'D45', 'E45', 'B46', 'C46', 'D46', 'E46', 'B47', 'C47', 'D47', 'E47', 'B48', 'C48', 'D48', 'E48', 'B49', 'C49', 'D49', 'E49', 'B50', 'C50', 'D50', 'E50', 'B51', 'C51', 'D51', 'E51')
foreach ($c in $cells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range('D2').Value = '30.350.87'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.862.41'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '233.67'
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('D7').Value = '0.4768'
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').Value = '0.2758'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('D9').Value = '0.06454'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').Value = '1.873.79'
$ws.Range('E10').Value = '  -6.20%  '
$ws.Range('D11').Value = '0.07434'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('D13').Value = '4.996'
$ws.Range('D14').Value = '85.88'
$ws.Range('E14').Value = '  -3.09%  '
$ws.Range('D15').Value = '0.6336'
$ws.Range('E15').Value = '  -4.30%  '
$ws.Range('D16').Value = '30.310.79'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').Value = '0.9998'
$ws.Range('D18').Value = '231.76'
$ws.Range('E18').Value = '  +2.39%  '
$ws.Range('D19').Value = '12.81'
$ws.Range('E19').Value = '  -3.82%  '
$ws.Range('D20').Value = '0.000007384'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').Value = '2.096.78'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '5.095'
$ws.Range('B24').Value = 'BitDAO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D24').Value = '0.3933'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = '6.015'
$ws.Range('E25').Value = '  -3.11%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '9.300'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '167.47'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '17.91'
$ws.Range('E28').Value = '  -4.04%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '1.864'
$ws.Range('E29').Value = '  -4.57%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.383'
$ws.Range('E30').Value = '  -5.41%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.09971'
$ws.Range('E31').Value = '  +5.05%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.222'
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '3.925'
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.04920'
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.150'
$ws.Range('E35').Value = '  -4.95%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7246'
$ws.Range('E36').Value = '  -3.25%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = '0.9993'
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.696'
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01940'
$ws.Range('E39').Value = '  +6.07%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.633'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '0.9048'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '1.985'
$ws.Range('E42').Value = '  -4.19%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '105.57'
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.4113'
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '5.546'
$ws.Range('E46').Value = '  -4.72%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '7.073'
$ws.Range('E47').Value = '  -5.37%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '61.34'
$ws.Range('E48').Value = '  -4.76%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.1209'
$ws.Range('E49').Value = '  -5.83%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.793'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.401'
$ws.Range('E51').Value = '  -4.93%  '
